$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new columns I ("I0") and J ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, borders, centered alignment) from the
# existing "IP" header cell (H1) onto the two new header cells so the
# style matches the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows 2-25: column I is always 1, column J duplicates column H ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
